# "Processed Salesforce Current Translations"
#
# The translation-validation data (rows for "Additional Info", "Applicable",
# and the shipment-identifier error message) was previously sitting on the
# "Text Not Translated" sheet with "Misssing" placeholders in the
# package/packageName/masterInfoLabel columns (O:Q). Salesforce processing
# has now resolved those placeholders to real values, and the fully-resolved
# table moves to the "Text Translated" sheet, leaving "Text Not Translated"
# empty again.

$wb = $excel.ActiveWorkbook

$wsTranslated = $wb.Worksheets.Item("Text Translated")
$wsNotTranslated = $wb.Worksheets.Item("Text Not Translated")

# Move (cut+paste) the whole A1:U4 table from "Text Not Translated" to
# "Text Translated" - this also clears the source sheet back to empty.
$wsNotTranslated.Range("A1:U4").Cut($wsTranslated.Range("A1"))

# Fill in the previously-"Misssing" package / packageName / masterInfoLabel
# columns now that the Salesforce translation metadata has been resolved.

# Row 2: "Additional Info"
$wsTranslated.Range("O2").Value = 'CustomLabel$CEC_Shipment_Additional_Info'
$wsTranslated.Range("P2").Value = 'CustomLabel'
$wsTranslated.Range("Q2").Value = 'Additional Info'

# Row 3: "Applicable"
$wsTranslated.Range("O3").Value = 'CustomField$CEC_ShippingIdentifier__c.CEC_Applicable__c | CustomLabel$CEC_ShipmentIdentifier_Applicable'
$wsTranslated.Range("P3").Value = 'CustomField | CustomLabel'
$wsTranslated.Range("Q3").Value = 'Applicable'

# Row 4: shipment-identifier error message
$wsTranslated.Range("O4").Value = 'CustomLabel$CEC_Shipment_Non_Applicable_Error'
$wsTranslated.Range("P4").Value = 'CustomLabel'
$wsTranslated.Range("Q4").Value = 'You must select one or more shipment identifiers in order to proceed.'
